$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> ECs ---
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.455159666666667
$ws.Range("H2").Value = 10.365479
$ws.Range("I2").Value = 0.557269825537176
$ws.Range("J2").Value = 0.557269825537176
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 6.220104333333333
$ws.Range("N2").Value = 18.660313
$ws.Range("O2").Value = 0.9231675752574263
$ws.Range("P2").Value = 0.9231675752574263
$ws.Range("Q2").Value = 21.49145361499189
$ws.Range("R2").Value = 193.423082534927
$ws.Range("S2").Value = 0.5144534336052837
$ws.Range("T2").Value = 0.5144534336052837

# --- Row 3: ECs -> FAPs ---
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.455159666666667
$ws.Range("H3").Value = 10.365479
$ws.Range("I3").Value = 0.557269825537176
$ws.Range("J3").Value = 0.557269825537176
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.09168666666666665
$ws.Range("N3").Value = 0.27506
$ws.Range("O3").Value = 0.01360783569119702
$ws.Range("P3").Value = 0.01360783569119702
$ws.Range("Q3").Value = 0.3167920726377778
$ws.Range("R3").Value = 2.85112865374
$ws.Range("S3").Value = 0.007583236221571918
$ws.Range("T3").Value = 0.007583236221571918

# --- Row 4: ECs -> sCs ---
$ws.Range("A4").Value = "ECs"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.455159666666667
$ws.Range("H4").Value = 10.365479
$ws.Range("I4").Value = 0.557269825537176
$ws.Range("J4").Value = 0.557269825537176
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4259936666666667
$ws.Range("N4").Value = 1.277981
$ws.Range("O4").Value = 0.06322458905137664
$ws.Range("P4").Value = 0.06322458905137664
$ws.Range("Q4").Value = 1.471876135322111
$ws.Range("R4").Value = 13.246885217899
$ws.Range("S4").Value = 0.0352331557103203
$ws.Range("T4").Value = 0.0352331557103203

# --- Row 5: FAPs -> ECs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Wnt4"
$ws.Range("C5").Value = "Fzd6"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.884975666666667
$ws.Range("H5").Value = 5.654927
$ws.Range("I5").Value = 0.3040207001254323
$ws.Range("J5").Value = 0.3040207001254323
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 6.220104333333333
$ws.Range("N5").Value = 18.660313
$ws.Range("O5").Value = 0.9231675752574263
$ws.Range("P5").Value = 0.9231675752574263
$ws.Range("Q5").Value = 11.72474531246122
$ws.Range("R5").Value = 105.522707812151
$ws.Range("S5").Value = 0.2806620525628605
$ws.Range("T5").Value = 0.2806620525628605

# --- Row 6: FAPs -> FAPs ---
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Wnt4"
$ws.Range("C6").Value = "Fzd6"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.884975666666667
$ws.Range("H6").Value = 5.654927
$ws.Range("I6").Value = 0.3040207001254323
$ws.Range("J6").Value = 0.3040207001254323
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.09168666666666665
$ws.Range("N6").Value = 0.27506
$ws.Range("O6").Value = 0.01360783569119702
$ws.Range("P6").Value = 0.01360783569119702
$ws.Range("Q6").Value = 0.1728271356244444
$ws.Range("R6").Value = 1.55544422062
$ws.Range("S6").Value = 0.004137063734029563
$ws.Range("T6").Value = 0.004137063734029563

# --- Row 7: FAPs -> sCs ---
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Wnt4"
$ws.Range("C7").Value = "Fzd6"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.884975666666667
$ws.Range("H7").Value = 5.654927
$ws.Range("I7").Value = 0.3040207001254323
$ws.Range("J7").Value = 0.3040207001254323
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4259936666666667
$ws.Range("N7").Value = 1.277981
$ws.Range("O7").Value = 0.06322458905137664
$ws.Range("P7").Value = 0.06322458905137664
$ws.Range("Q7").Value = 0.8029876958207778
$ws.Range("R7").Value = 7.226889262387
$ws.Range("S7").Value = 0.01922158382854227
$ws.Range("T7").Value = 0.01922158382854227

# --- Row 8: sCs -> ECs ---
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Wnt4"
$ws.Range("C8").Value = "Fzd6"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8600203333333334
$ws.Range("H8").Value = 2.580061
$ws.Range("I8").Value = 0.1387094743373916
$ws.Range("J8").Value = 0.1387094743373916
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 6.220104333333333
$ws.Range("N8").Value = 18.660313
$ws.Range("O8").Value = 0.9231675752574263
$ws.Range("P8").Value = 0.9231675752574263
$ws.Range("Q8").Value = 5.349416202121444
$ws.Range("R8").Value = 48.144745819093
$ws.Range("S8").Value = 0.128052089089282
$ws.Range("T8").Value = 0.128052089089282

# --- Row 9: sCs -> FAPs ---
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Wnt4"
$ws.Range("C9").Value = "Fzd6"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8600203333333334
$ws.Range("H9").Value = 2.580061
$ws.Range("I9").Value = 0.1387094743373916
$ws.Range("J9").Value = 0.1387094743373916
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.09168666666666665
$ws.Range("N9").Value = 0.27506
$ws.Range("O9").Value = 0.01360783569119702
$ws.Range("P9").Value = 0.01360783569119702
$ws.Range("Q9").Value = 0.07885239762888888
$ws.Range("R9").Value = 0.70967157866
$ws.Range("S9").Value = 0.001887535735595535
$ws.Range("T9").Value = 0.001887535735595535

# --- Row 10: sCs -> sCs ---
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Wnt4"
$ws.Range("C10").Value = "Fzd6"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8600203333333334
$ws.Range("H10").Value = 2.580061
$ws.Range("I10").Value = 0.1387094743373916
$ws.Range("J10").Value = 0.1387094743373916
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4259936666666667
$ws.Range("N10").Value = 1.277981
$ws.Range("O10").Value = 0.06322458905137664
$ws.Range("P10").Value = 0.06322458905137664
$ws.Range("Q10").Value = 0.3663632152045556
$ws.Range("R10").Value = 3.297268936841
$ws.Range("S10").Value = 0.00876984951251406
$ws.Range("T10").Value = 0.00876984951251406

Write-Host "done"
